# Apply "updated mock data files" edit to the Negative Comments mock data workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data set (header + 20 comment rows), columns A:C.
$data = @(
    @("CLINIC", "RESPONSE", "COMMENTS"),
    @("A&E", "Extremely Unlikely", "Waited too long to find a parking spot"),
    @("Rehab Services", "Extremely Unlikely", "staff was rude"),
    @("Radiology", "Extremely Unlikely", "Food was terrible"),
    @("Day Surgery", "Unlikely", "Service recieved was adaquete but staff seemed like they dont care at all"),
    @("A&E", "Extremely Unlikely", "Waited for long time for poor service"),
    @("Special Care Baby Unit", "Extremely Unlikely", "Waited over 5 hours"),
    @("Day Surgery", "Extremely Unlikely", "Food was terrible"),
    @("Special Care Baby Unit", "Extremely Unlikely", "Felt as if i was not a priority"),
    @("Rehab Services", "Extremely Unlikely", "Clenliness isn't the best but otherwise okay"),
    @("Special Care Baby Unit", "Extremely Unlikely", "Waited for long time for poor service"),
    @("Orthodontics", "Unlikely", "Service recieved was adaquete but staff seemed like they dont care at all"),
    @("Sitwell", "Unlikely", "Waited too long to find a parking spot"),
    @("A&E", "Unlikely", "Waited over 5 hours"),
    @("Heart Failure", "Extremely Unlikely", "Food was terrible"),
    @("Bone Health", "Unlikely", "Doctors are patronising and made me feel bad"),
    @("Theatre Treatment Suite Implants", "Unlikely", "Long wait times"),
    @("A&E", "Extremely Unlikely", "Long wait times"),
    @("Radiology", "Extremely Unlikely", "doctors dont seem to care about me, felt ignored"),
    @("Bone Health", "Extremely Unlikely", "staff was rude"),
    @("Bone Health", "Extremely Unlikely", "staff was rude")
)

# Clear any previous contents below the new data range (workbook only had 6 rows before).
$ws.Cells.Clear()

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 1
    $rowVals = $data[$i]
    $ws.Cells.Item($rowNum, 1).Value = $rowVals[0]
    $ws.Cells.Item($rowNum, 2).Value = $rowVals[1]
    $ws.Cells.Item($rowNum, 3).Value = $rowVals[2]
}

$ws.Range("F21").Select()
